$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are stored as literal text in the sheet (prices/percentages as
# strings), so each one is written with a leading apostrophe to force a
# text entry instead of Excel auto-converting to a Number/Percentage.
$updates = @{
    "D2"  = "'324.29";    "E2"  = "'9.11%"
    "D3"  = "'49.71";     "E3"  = "'18.95%"
    "D4"  = "'5.336";     "E4"  = "'6.70%"
    "D5"  = "'0.08154";   "E5"  = "'8.49%"
    "E6"  = "'5.44%"
    "D7"  = "'1.664";     "E7"  = "'5.15%"
    "D8"  = "'1.168";     "E8"  = "'26.08%"
    "D9"  = "'0.1357";    "E9"  = "'13.35%"
    "D10" = "'0.1954";    "E10" = "'7.28%"
    "D11" = "'0.09535";   "E11" = "'6.58%"
    "D12" = "'0.04549";   "E12" = "'11.40%"
    "E13" = "'-0.05%"
    "D14" = "'0.001333";  "E14" = "'4.16%"
    "D15" = "'0.005944";  "E15" = "'2.28%"
    "D16" = "'3.395";     "E16" = "'1.13%"
    "E17" = "'1.54%"
    "D18" = "'0.3394";    "E18" = "'2.43%"
    "D19" = "'8.182";     "E19" = "'0.81%"
    "E20" = "'2.09%"
    "D21" = "'0.3053";    "E21" = "'-1.57%"
    "D22" = "'0.04300";   "E22" = "'5.03%"
    "D23" = "'0.001307";  "E23" = "'3.32%"
    "D24" = "'0.004263";  "E24" = "'9.40%"
    "E25" = "'9.62%"
    "D26" = "'0.0003722"; "E26" = "'-0.05%"
    "D38" = "'0.02783";   "E38" = "'15.77%"
    "D39" = "'0.05559";   "E39" = "'6.81%"
    "D40" = "'0.006301";  "E40" = "'-0.08%"
    "D41" = "'0.007690";  "E41" = "'-1.51%"
    "D42" = "'0.1449";    "E42" = "'9.21%"
    "D43" = "'0.007696";  "E43" = "'3.89%"
    "D44" = "'0.008093";  "E44" = "'11.33%"
    "D45" = "'0.3524";    "E45" = "'18.88%"
    "D46" = "'0.00006777";"E46" = "'2.85%"
    "E47" = "'-0.08%"
    "E48" = "'94.20%"
    "D49" = "'0.004000";  "E49" = "'-4.84%"
    "E50" = "'-0.08%"
    "E51" = "'-0.08%"
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
